$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New row 7: "07_20210603" / increase lakebed leakance / previously calibrated
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = "07_20210603"
$ws.Range("B7").Value = 10.199999999999999
$ws.Range("C7").NumberFormat = "0.00E+00"
$ws.Range("C7").Value = 0.00001
$ws.Range("D7").Value = 8.4
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 200
$ws.Range("G7").Value = 9.1999999999999993
$ws.Range("H7").NumberFormat = "0.00E+00"
$ws.Range("H7").Value = 0.00001
$ws.Range("I7").Value = 9.6
$ws.Range("J7").WrapText = $true
$ws.Range("J7").Value = "increase lakebed leakance"
$ws.Range("K7").WrapText = $true
$ws.Range("K7").Value = "previously calibrated"

# ---------------------------------------------------------------------------
# New row 8: "08_20210603" / increase lakebed leakance / increase Ks for zone
# containing H0_23
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "08_20210603"
$ws.Range("B8").Value = 10.199999999999999
$ws.Range("C8").NumberFormat = "0.00E+00"
$ws.Range("C8").Value = 0.00001
$ws.Range("D8").Value = 8.4
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 200
$ws.Range("G8").Value = 9.1999999999999993
$ws.Range("H8").NumberFormat = "0.00E+00"
$ws.Range("H8").Value = 0.00001
$ws.Range("I8").Value = 9.6
$ws.Range("J8").WrapText = $true
$ws.Range("J8").Value = "increase lakebed leakance"

# ---------------------------------------------------------------------------
# Column K header text changed from "Vertical K" to "Ks"
# ---------------------------------------------------------------------------
$ws.Range("K1").Value = "Ks"

$ws.Range("K8").WrapText = $true
$ws.Range("K8").Value = "increase Ks  for zone containing H0_23"

# ---------------------------------------------------------------------------
# Row heights for the two new rows (match wrapped-text autofit heights)
# ---------------------------------------------------------------------------
$ws.Rows.Item(7).RowHeight = 45
$ws.Rows.Item(8).RowHeight = 60

# ---------------------------------------------------------------------------
# Selection moved to K8
# ---------------------------------------------------------------------------
$null = $ws.Range("K8").Select()
